# edit.ps1 - Applies the "cosmic patterns" -> "wonders of science" rewrite
# described by the commit diff, using the Word COM-interop object model.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title (Paragraph 1)
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1.SetRange($p1.Start, $p1.End - 1)
$p1.Text = "Unlocking the Wonders of Science: A Journey Through Mathematics, Chemistry, Biology, and Medicine"

# ---------------------------------------------------------------------------
# 2. Author name (Paragraph 2): "David Brooks" -> "Dr. Evelyn Coleman"
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2).Range
$p2.SetRange($p2.Start, $p2.End - 1)
$p2.Text = "Dr. Evelyn Coleman"

# ---------------------------------------------------------------------------
# 3. Email (Paragraph 3): "davidbrooks@hotmail.com" -> "evelyn.coleman10@highereducation.edu"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$p3.SetRange($p3.Start, $p3.End - 1)
$p3.Text = "evelyn.coleman10@highereducation.edu"

# ---------------------------------------------------------------------------
# 4. Main body (Paragraph 5): full rewrite from "cosmic patterns" essay to the
#    "wonders of science" essay, including the manual line breaks (vertical
#    tab, Chr(11)) that separate each mini-section.
# ---------------------------------------------------------------------------
$NL = [char]11
$bodyText = "Science, a realm that unveils the mysteries of the universe, invites us on an awe-inspiring voyage of discovery." +
    " From the intricate beauty of mathematical patterns to the wonders of chemical reactions, from the intricate workings of living organisms to the marvels of medicine, science offers a plethora of knowledge and endless opportunities for exploration" +
    $NL + $NL +
    "Mathematics, the universal language of science, provides the tools to decipher complex concepts and quantify the world around us." +
    " Through equations, formulas, and algorithms, we unravel the secrets of numbers and shapes." +
    " From simple arithmetic to calculus and beyond, mathematics serves as a cornerstone for advancements in various disciplines." +
    $NL + $NL +
    "Chemistry, the study of matter and its interactions, delves into the molecular realm, revealing the fundamental building blocks of the universe." +
    " By exploring atomic structure, chemical bonding, and reactions, we gain insights into the diverse properties of substances." +
    " From everyday materials to pharmaceuticals, chemistry touches every aspect of our lives." +
    $NL + $NL +
    "Introduction Continued:" +
    $NL + $NL +
    "Biology, the science of life, investigates the intricate complexities of living organisms." +
    " From the microscopic world of cells to the diversity of ecosystems, we delve into the mysteries ofSheng Ming De  origin, evolution, and behavior." +
    " By understanding the intricacies of life, we gain valuable insights into preserving and promoting health and well-being." +
    $NL + $NL +
    "Medicine, a practical application of scientific knowledge, diagnoses, treats, and prevents diseases." +
    " By harnessing advances in chemistry, biology, and other disciplines, medical researchers and practitioners strive to improve human health and alleviate suffering." +
    " From antibiotics to vaccines, medicine has transformed countless lives and continues to hold the promise of eradicating diseases." +
    $NL + $NL +
    "Introduction Concluded:" +
    $NL + $NL +
    "Science is an ever-evolving tapestry, continuously enriched by new discoveries and insights." +
    " It empowers us to understand the world around us, to solve complex problems, and to create innovative technologies that improve our lives." +
    " Embarking on this journey of scientific exploration promises a rewarding and fulfilling experience, fostering curiosity, critical thinking, and a deeper appreciation for the beauty of the natural world."

$p5 = $d.Paragraphs(5).Range
$p5.SetRange($p5.Start, $p5.End - 1)
$p5.Text = $bodyText

# ---------------------------------------------------------------------------
# 5. Summary body (Paragraph 7): rewrite closing summary paragraph.
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7).Range
$p7.SetRange($p7.Start, $p7.End - 1)
$p7.Text = "This essay takes readers on an enthralling journey through the interconnected world of science, exploring the depths of mathematics, chemistry, biology, and medicine." +
    " From unraveling mathematical patterns to deciphering chemical reactions, from understanding the intricacies of life to harnessing scientific knowledge for medical advancements, science unveils the wonders of the universe and offers endless opportunities for exploration and discovery."

# ---------------------------------------------------------------------------
# 6. Add a trailing empty paragraph at the very end of the document body.
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
